# Minor text updates to two slides of "02 - Structure of Compilers.pptx"
$p = $ppt.ActivePresentation

# --- Slide 18 ("Passes"): clarify that the note refers to a compiler "pass" ---
$slide18 = $p.Slides.Item(18)
$body18  = $slide18.Shapes.Item(4).TextFrame.TextRange
$notePara = $body18.Paragraphs(3)
# Replace the single run's text in place so the run/paragraph formatting is preserved.
$notePara.Runs(1).Text = "Note: Some authors restrict the definition of compiler pass to a traversal that involves disk I/O, but we will use a more general definition."

# --- Slide 19 ("Single-pass Versus Multi-pass Compilers"): reword the "multiprocessor" bullet ---
$slide19 = $p.Slides.Item(19)
$body19  = $slide19.Shapes.Item(4).TextFrame.TextRange
$multiPara = $body19.Paragraphs(5)
# Rewrite the start of the bullet, then type the remainder of the new phrase
# after it, so the text ends up split across two runs (matching how the
# sentence was actually extended during editing).
$multiPara.Text = "can exploit concurrency and "
$multiPara.InsertAfter("multiprocessor architectures") | Out-Null
